$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@(2, 338639, 430791589)
    ,@(3, 276, 331062)
    ,@(8, 896, 1317211)
    ,@(10, 122138, 178889619)
    ,@(12, 63262, 91297391)
    ,@(16, 4098, 5817021)
    ,@(21, 7870, 10988003)
    ,@(23, 81157, 100919336)
    ,@(29, 33626, 49209492)
    ,@(32, 12091, 17391923)
    ,@(37, 2130, 3008885)
    ,@(38, 101477, 127393573)
    ,@(40, 83, 116508)
    ,@(46, 45923, 67274427)
    ,@(48, 9632, 13802163)
    ,@(53, 2711, 3795439)
    ,@(54, 72345, 90595449)
    ,@(61, 29321, 42986372)
    ,@(64, 11750, 16981173)
    ,@(70, 1725, 2419319)
    ,@(72, 21640, 28314101)
    ,@(76, 8025, 11751804)
    ,@(78, 5433, 7895983)
    ,@(79, 526, 746397)
    ,@(80, 323, 456323)
    ,@(81, 147915, 184032906)
    ,@(87, 66021, 96723033)
    ,@(90, 31192, 45110249)
    ,@(92, 2883, 4155198)
    ,@(94, 3270, 4614362)
    ,@(95, 36075, 48918778)
    ,@(99, 8794, 12922758)
    ,@(101, 8182, 11869214)
    ,@(103, 578, 820376)
    ,@(105, 13831, 23978665)
    ,@(108, 3275, 6116377)
    ,@(110, 4575, 8579021)
    ,@(112, 198, 363790)
    ,@(113, 279, 503191)
    ,@(115, 148702, 183731493)
    ,@(121, 54971, 80514899)
    ,@(123, 29010, 42018787)
    ,@(124, 1346, 1843694)
    ,@(127, 2606, 3667952)
    ,@(129, 569252, 751478166)
    ,@(134, 1479, 2190317)
    ,@(135, 33, 43510)
    ,@(136, 224130, 329311239)
    ,@(139, 203768, 296259699)
    ,@(142, 2946, 4136972)
    ,@(145, 7588, 10704714)
    ,@(147, 17, 25500)
    ,@(148, 47751, 63669096)
    ,@(154, 14801, 21690198)
    ,@(155, 3999, 5770075)
    ,@(160, 452, 639898)
    ,@(161, 18823, 24888603)
    ,@(165, 7737, 11257767)
    ,@(167, 5451, 7852391)
    ,@(172, 24377, 45524979)
    ,@(173, 2495, 4608491)
    ,@(178, 91806, 114516009)
    ,@(185, 35270, 51697389)
    ,@(187, 13806, 19944343)
    ,@(189, 1287, 1801399)
    ,@(191, 1889, 2652271)
    ,@(193, 249923, 310081612)
    ,@(195, 179, 257979)
    ,@(199, 916, 1346590)
    ,@(201, 89909, 131746654)
    ,@(204, 34880, 50213587)
    ,@(207, 5298, 7546779)
    ,@(209, 14, 21000)
    ,@(210, 5573, 7723766)
    ,@(213, 277301, 342757810)
    ,@(214, 166, 185413)
    ,@(220, 639, 929514)
    ,@(222, 99296, 145223687)
    ,@(223, 101, 150199)
    ,@(225, 54508, 78749310)
    ,@(228, 4829, 6778453)
    ,@(231, 6848, 9487346)
    ,@(234, 111812, 139544393)
    ,@(236, 82, 115713)
    ,@(239, 582, 850939)
    ,@(241, 51463, 75376843)
    ,@(243, 13310, 19142406)
    ,@(245, 1933, 2771963)
    ,@(247, 2863, 4011238)
    ,@(248, 272282, 343520374)
    ,@(254, 870, 1277868)
    ,@(256, 100303, 146938850)
    ,@(259, 69396, 100580902)
    ,@(261, 2496, 3518644)
    ,@(264, 5372, 7537580)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 3).Value = $u[1]
    $ws.Cells.Item($row, 4).Value = $u[2]
}

$wb.Save()
